$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table updates (rows 3-5) ---
$ws.Range("C3").Value = 8919
$ws.Range("D3").Value = 35.9

$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 6657
$ws.Range("D4").Value = 94.40000000000001

$ws.Range("B5").Value = 11
$ws.Range("C5").Value = 15576

# --- "Good Drivers" table: a new driver was added at the top of the list ---
# Insert a new row at row 13, shifting the existing rows (13-18) down to (14-19).
$ws.Rows.Item(13).Insert()

# Fill in the newly inserted row with the new driver's data.
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B13").Value = 11128
$ws.Range("B13").NumberFormat = "#,##0"
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = ""

# Update the client counts for the rest of the (shifted) good-driver rows.
$ws.Range("B14").Value = 486214
$ws.Range("B15").Value = 79953
$ws.Range("B16").Value = 35355
$ws.Range("B17").Value = 65425
$ws.Range("B18").Value = 117653
# Row 19 (Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1) keeps its original values.
